# Auto-generated edit script: refresh the cryptocurrency Price (D) and
# Volume(1h) (E) columns on the active sheet to the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# cell -> new value
$updates = @{
    'D2' = '26.594.42'
    'E2' = '  +0.00%  '
    'D3' = '1.595.18'
    'E3' = '  +0.38%  '
    'E4' = '  -0.01%  '
    'D5' = '210.93'
    'E5' = '  +0.14%  '
    'E6' = '  +1.11%  '
    'E7' = '  -0.05%  '
    'E8' = '  +0.12%  '
    'E9' = '  -0.97%  '
    'D10' = '19.43'
    'E10' = '  -0.57%  '
    'E11' = '  +0.06%  '
    'D12' = '1.819.65'
    'E12' = '  +0.41%  '
    'D13' = '1.574.12'
    'E13' = '  -0.99%  '
    'D14' = '4.01'
    'E14' = '  -0.07%  '
    'E15' = '  -0.33%  '
    'D16' = '64.35'
    'E16' = '  -0.41%  '
    'D17' = '26.594.19'
    'E17' = '  -0.04%  '
    'D18' = '0.0₃0730'
    'E18' = '  +0.62%  '
    'E19' = '  +0.00%  '
    'D20' = '207.80'
    'E20' = '  -0.18%  '
    'D21' = '6.89'
    'D22' = '4.24'
    'E22' = '  +0.12%  '
    'E23' = '  -1.76%  '
    'D24' = '8.84'
    'E24' = '  -0.03%  '
    'E25' = '  -0.93%  '
    'E26' = '  -0.11%  '
    'D27' = '7.12'
    'E27' = '  -1.55%  '
    'E28' = '  +0.47%  '
    'D29' = '15.21'
    'E29' = '  -0.35%  '
    'D30' = '0.0504'
    'E30' = '  -0.62%  '
    'E31' = '  +0.32%  '
    'E32' = '  +0.13%  '
    'D33' = '0.653'
    'E33' = '  -1.69%  '
    'D34' = '2.91'
    'E34' = '  +0.67%  '
    'D35' = '1.280.42'
    'E35' = '  -2.01%  '
    'E36' = '  +0.83%  '
    'E37' = '  +0.12%  '
    'D39' = '0.842'
    'E39' = '  +1.70%  '
    'E40' = '  +0.03%  '
    'E41' = '  +1.71%  '
    'E42' = '  +1.45%  '
    'D43' = '0.785'
    'E43' = '  -0.54%  '
    'D44' = '63.72'
    'E44' = '  +1.55%  '
    'E45' = '  +9.51%  '
    'D46' = '1.731.81'
    'E46' = '  +0.38%  '
    'D47' = '89.47'
    'E47' = '  -0.33%  '
    'E48' = '  -1.15%  '
    'E49' = '  -0.39%  '
    'E50' = '  +4.56%  '
    'E51' = '  +0.85%  '
}

# Columns D/E hold plain text in the source sheet (prices use "."
# as a thousands separator and percentages keep surrounding padding),
# so digit-only values are entered with a leading apostrophe to force
# text storage instead of being auto-converted to a number, then the
# style is reset to "Normal" so no stray number-format gets attached.
foreach ($cellRef in $updates.Keys) {
    $value = $updates[$cellRef]
    $range = $ws.Range($cellRef)
    $isNumericLooking = $value -match "^[+-]?[0-9]*\.?[0-9]+$"
    if ($isNumericLooking) {
        $range.Value = "'" + $value
        $range.Style = 'Normal'
    } else {
        $range.Value = $value
    }
}
